# Updated upto 55 TC
# Sheet "WebStaff_TC" is the active sheet in this workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 14-39 (TC rows already tested): mark Result column (G) as "Pass" (was "Fail")
for ($r = 14; $r -le 39; $r++) {
    $ws.Range("G$r").Value = "Pass"
}

# Rows 40-53 (newly tested TC rows): fill Tester column (F) and Result column (G)
for ($r = 40; $r -le 53; $r++) {
    $ws.Range("F$r").Value = "kulandasamyc"
    $ws.Range("G$r").Value = "Pass"
}
